{"js": "// The scraped diff shows the document's styles part gaining a second,\n// redundant copy of the built-in style set (Normal, Table Normal,\n// Heading 1-6, Title, Subtitle) immediately before the pre-existing\n// \"Subtitle\" style, plus two new Google-Docs round-trip custom XML\n// parts (customXML/item1.xml + itemProps1.xml) tacked on at package\n// level. Nothing in the body text changes.\n//\n// That shape is a save-time artifact of whatever produced the \"after\"\n// package (re-serializing its own built-in style catalog on top of the\n// one already present) rather than a user-visible formatting edit, so\n// there is no dedicated \"duplicate this style\" verb in the Word\n// JavaScript API - style identities are unique by design. The closest,\n// faithful expression of the same intent through the supported object\n// model is to (re)assert every one of those built-in styles so the\n// document is guaranteed to carry a fully defined copy of each one,\n// and to register the same custom XML payloads the diff adds.\n\n// 1) Make sure every built-in style referenced by the diff is present /\n//    re-applied. For styles Word already ships (Normal, Heading 1-6,\n//    Title, Subtitle, Table Normal) this is idempotent against the\n//    existing definition - exactly the no-visual-effect \"rewrite\" the\n//    diff captures.\nconst builtInStyles = [\n  { name: \"Normal\", type: Word.StyleType.paragraph },\n  { name: \"Table Normal\", type: Word.StyleType.table },\n  { name: \"Heading 1\", type: Word.StyleType.paragraph },\n  { name: \"Heading 2\", type: Word.StyleType.paragraph },\n  { name: \"Heading 3\", type: Word.StyleType.paragraph },\n  { name: \"Heading 4\", type: Word.StyleType.paragraph },\n  { name: \"Heading 5\", type: Word.StyleType.paragraph },\n  { name: \"Heading 6\", type: Word.StyleType.paragraph },\n  { name: \"Title\", type: Word.StyleType.paragraph },\n  { name: \"Subtitle\", type: Word.StyleType.paragraph },\n];\n\nfor (const s of builtInStyles) {\n  try {\n    context.document.addStyle(s.name, s.type);\n  } catch (e) {\n    // Some hosts throw if a style with that name is already defined;\n    // that's fine, it means the style is already present.\n  }\n}\nawait context.sync();\n\n// 2) Mirror the two custom XML parts the diff adds\n//    (customXML/item1.xml + customXML/itemProps1.xml roundtrip data).\ntry {\n  const customXmlParts = context.document.customXmlParts;\n  customXmlParts.add(\n    '<?xml version=\"1.0\" encoding=\"utf-8\"?>' +\n      '<go:gDocsCustomXmlDataStorage xmlns:go=\"http://customooxmlschemas.google.com/\" ' +\n      'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">' +\n      '<go:docsCustomData xmlns:go=\"http://customooxmlschemas.google.com/\" ' +\n      'roundtripDataSignature=\"AMtx7miKB8iE5r8vC0oZ8nOjxRozQS8wqA==\">' +\n      \"AMUW2mXriySWwlJ5Eg4iFDWJHSA1NHyNoWbA4tmYEDTfwkY0gBAmHyyseNjoyh3G/EqUH5nUN0Mw3/zACVp7JxHh0cFsSYm+cd0bm4VJgXd2ytoUjP9Vi+U=\" +\n      \"</go:docsCustomData></go:gDocsCustomXmlDataStorage>\"\n  );\n  await context.sync();\n} catch (e) {\n  // Custom XML parts are an optional, host-dependent capability -\n  // ignore if unsupported rather than fail the whole edit.\n}\n", "ps1": "# The scraped diff shows the document's styles part gaining a second,\n# redundant copy of the built-in style set (Normal, Table Normal,\n# Heading 1-6, Title, Subtitle) immediately before the pre-existing\n# \"Subtitle\" style, plus two new Google-Docs round-trip custom XML\n# parts (customXML/item1.xml + itemProps1.xml) tacked on at package\n# level. Nothing in the body text changes.\n#\n# That shape is a save-time artifact of whatever produced the \"after\"\n# package (re-serializing its own built-in style catalog on top of the\n# one already present) rather than a user-visible formatting edit, so\n# there is no dedicated \"duplicate this style\" verb in the Word object\n# model - style identities are unique by design (Styles.Add on an\n# existing name just hands back the existing Style). The closest,\n# faithful expression of the same intent through the supported object\n# model is to (re)assert every one of those built-in styles so the\n# document is guaranteed to carry a fully defined copy of each one, and\n# to register the same custom XML payloads the diff adds.\n\n$d = $word.ActiveDocument\n\n# 1) Make sure every built-in style referenced by the diff is present /\n#    re-applied. For styles Word already ships (Normal, Heading 1-6,\n#    Title, Subtitle, Table Normal) this is idempotent against the\n#    existing definition - exactly the no-visual-effect \"rewrite\" the\n#    diff captures.\n$wdStyleTypeParagraph = 1\n$wdStyleTypeTable = 3\n\n$builtInStyles = @(\n    @{ Name = \"Normal\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Table Normal\"; Type = $wdStyleTypeTable },\n    @{ Name = \"Heading 1\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Heading 2\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Heading 3\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Heading 4\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Heading 5\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Heading 6\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Title\"; Type = $wdStyleTypeParagraph },\n    @{ Name = \"Subtitle\"; Type = $wdStyleTypeParagraph }\n)\n\nforeach ($s in $builtInStyles) {\n    try {\n        $d.Styles.Add($s.Name, $s.Type) | Out-Null\n    } catch {\n        # Already defined with that name - the style is present either way.\n    }\n}\n\n# 2) Mirror the two custom XML parts the diff adds\n#    (customXML/item1.xml + customXML/itemProps1.xml roundtrip data).\ntry {\n    $customXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>' +\n        '<go:gDocsCustomXmlDataStorage xmlns:go=\"http://customooxmlschemas.google.com/\" ' +\n        'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">' +\n        '<go:docsCustomData xmlns:go=\"http://customooxmlschemas.google.com/\" ' +\n        'roundtripDataSignature=\"AMtx7miKB8iE5r8vC0oZ8nOjxRozQS8wqA==\">' +\n        'AMUW2mXriySWwlJ5Eg4iFDWJHSA1NHyNoWbA4tmYEDTfwkY0gBAmHyyseNjoyh3G/EqUH5nUN0Mw3/zACVp7JxHh0cFsSYm+cd0bm4VJgXd2ytoUjP9Vi+U=' +\n        '</go:docsCustomData></go:gDocsCustomXmlDataStorage>'\n    $d.CustomXMLParts.Add($customXml) | Out-Null\n} catch {\n    # Custom XML parts are an optional, host-dependent capability -\n    # ignore if unsupported rather than fail the whole edit.\n}\n\n$d.Save()\n"}
